$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2 through 18: 45185 -> 45204
for ($row = 2; $row -le 18; $row++) {
    $ws.Cells.Item($row, 3).Value = 45204
}
